$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (bold/border/centered style used by col A) down into the
# two newly added rows before filling their values, by cloning the format
# already present on A30.
$ws.Range("A30").Copy()
$ws.Range("A31:A32").PasteSpecial(-4122)

# Column A index values for newly added rows 31 and 32
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(32, 1).Value = 30

# Row 2
$ws.Cells.Item(2, 2).Value = "NSE:ABB"
$ws.Cells.Item(2, 3).Value = "NSE:AARON"
$ws.Cells.Item(2, 4).Value = "NSE:CAMS"
$ws.Cells.Item(2, 5).Value = "NSE:INOXWIND"
$ws.Cells.Item(2, 6).Value = "NSE:ABB"

# Row 3
$ws.Cells.Item(3, 2).Value = "NSE:AKASH"
$ws.Cells.Item(3, 3).Value = "NSE:BCG"
$ws.Cells.Item(3, 4).Value = "NSE:CDSL"
$ws.Cells.Item(3, 5).Value = ""
$ws.Cells.Item(3, 6).Value = "NSE:BANKBARODA"

# Row 4
$ws.Cells.Item(4, 2).Value = "NSE:ASHOKA"
$ws.Cells.Item(4, 3).Value = "NSE:CORALFINAC"
$ws.Cells.Item(4, 4).Value = "NSE:IEX"
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = "NSE:CAMS"

# Row 5
$ws.Cells.Item(5, 2).Value = "NSE:BANKBARODA"
$ws.Cells.Item(5, 3).Value = "NSE:DIACABS"
$ws.Cells.Item(5, 4).Value = "NSE:PNBHOUSING"
$ws.Cells.Item(5, 5).Value = ""
$ws.Cells.Item(5, 6).Value = "NSE:CIPLA"

# Row 6
$ws.Cells.Item(6, 2).Value = "NSE:BSOFT"
$ws.Cells.Item(6, 3).Value = "NSE:GMMPFAUDLR"
$ws.Cells.Item(6, 4).Value = ""
$ws.Cells.Item(6, 5).Value = ""
$ws.Cells.Item(6, 6).Value = "NSE:COFORGE"

# Row 7
$ws.Cells.Item(7, 2).Value = "NSE:CAMS"
$ws.Cells.Item(7, 3).Value = "NSE:GODFRYPHLP"
$ws.Cells.Item(7, 4).Value = ""
$ws.Cells.Item(7, 5).Value = ""
$ws.Cells.Item(7, 6).Value = "NSE:LODHA"

# Row 8
$ws.Cells.Item(8, 2).Value = "NSE:CAPACITE"
$ws.Cells.Item(8, 3).Value = "NSE:HINDWAREAP"
$ws.Cells.Item(8, 4).Value = ""
$ws.Cells.Item(8, 5).Value = ""
$ws.Cells.Item(8, 6).Value = "NSE:MARUTI"

# Row 9
$ws.Cells.Item(9, 2).Value = "NSE:CIPLA"
$ws.Cells.Item(9, 3).Value = "NSE:LIKHITHA"
$ws.Cells.Item(9, 4).Value = ""
$ws.Cells.Item(9, 5).Value = ""
$ws.Cells.Item(9, 6).Value = "NSE:OBEROIRLTY"

# Row 10
$ws.Cells.Item(10, 2).Value = "NSE:COFORGE"
$ws.Cells.Item(10, 3).Value = "NSE:MSPL"
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(10, 5).Value = ""
$ws.Cells.Item(10, 6).Value = ""

# Row 11
$ws.Cells.Item(11, 2).Value = "NSE:GNA"
$ws.Cells.Item(11, 3).Value = "NSE:NURECA"
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = ""
$ws.Cells.Item(11, 6).Value = ""

# Row 12
$ws.Cells.Item(12, 2).Value = "NSE:HDFCLIQUID"
$ws.Cells.Item(12, 3).Value = "NSE:S&SPOWER"
$ws.Cells.Item(12, 4).Value = ""
$ws.Cells.Item(12, 5).Value = ""
$ws.Cells.Item(12, 6).Value = ""

# Row 13
$ws.Cells.Item(13, 2).Value = "NSE:ICDSLTD"
$ws.Cells.Item(13, 3).Value = ""
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 5).Value = ""
$ws.Cells.Item(13, 6).Value = ""

# Row 14
$ws.Cells.Item(14, 2).Value = "NSE:INOXGREEN"
$ws.Cells.Item(14, 3).Value = ""
$ws.Cells.Item(14, 4).Value = ""
$ws.Cells.Item(14, 5).Value = ""
$ws.Cells.Item(14, 6).Value = ""

# Row 15
$ws.Cells.Item(15, 2).Value = "NSE:IRMENERGY"
$ws.Cells.Item(15, 3).Value = ""
$ws.Cells.Item(15, 4).Value = ""
$ws.Cells.Item(15, 5).Value = ""
$ws.Cells.Item(15, 6).Value = ""

# Row 16
$ws.Cells.Item(16, 2).Value = "NSE:ITBEES"
$ws.Cells.Item(16, 3).Value = ""
$ws.Cells.Item(16, 4).Value = ""
$ws.Cells.Item(16, 5).Value = ""
$ws.Cells.Item(16, 6).Value = ""

# Row 17
$ws.Cells.Item(17, 2).Value = "NSE:KDDL"
$ws.Cells.Item(17, 3).Value = ""
$ws.Cells.Item(17, 4).Value = ""
$ws.Cells.Item(17, 5).Value = ""
$ws.Cells.Item(17, 6).Value = ""

# Row 18
$ws.Cells.Item(18, 2).Value = "NSE:LTTS"
$ws.Cells.Item(18, 3).Value = ""
$ws.Cells.Item(18, 4).Value = ""
$ws.Cells.Item(18, 5).Value = ""
$ws.Cells.Item(18, 6).Value = ""

# Row 19
$ws.Cells.Item(19, 2).Value = "NSE:MARUTI"
$ws.Cells.Item(19, 3).Value = ""
$ws.Cells.Item(19, 4).Value = ""
$ws.Cells.Item(19, 5).Value = ""
$ws.Cells.Item(19, 6).Value = ""

# Row 20
$ws.Cells.Item(20, 2).Value = "NSE:MOL"
$ws.Cells.Item(20, 3).Value = ""
$ws.Cells.Item(20, 4).Value = ""
$ws.Cells.Item(20, 5).Value = ""
$ws.Cells.Item(20, 6).Value = ""

# Row 21
$ws.Cells.Item(21, 2).Value = "NSE:MONIFTY500"
$ws.Cells.Item(21, 3).Value = ""
$ws.Cells.Item(21, 4).Value = ""
$ws.Cells.Item(21, 5).Value = ""
$ws.Cells.Item(21, 6).Value = ""

# Row 22
$ws.Cells.Item(22, 2).Value = "NSE:NIF100BEES"
$ws.Cells.Item(22, 3).Value = ""
$ws.Cells.Item(22, 4).Value = ""
$ws.Cells.Item(22, 5).Value = ""
$ws.Cells.Item(22, 6).Value = ""

# Row 23
$ws.Cells.Item(23, 2).Value = "NSE:NIFTYQLITY"
$ws.Cells.Item(23, 3).Value = ""
$ws.Cells.Item(23, 4).Value = ""
$ws.Cells.Item(23, 5).Value = ""
$ws.Cells.Item(23, 6).Value = ""

# Row 24
$ws.Cells.Item(24, 2).Value = "NSE:NIITLTD"
$ws.Cells.Item(24, 3).Value = ""
$ws.Cells.Item(24, 4).Value = ""
$ws.Cells.Item(24, 5).Value = ""
$ws.Cells.Item(24, 6).Value = ""

# Row 25
$ws.Cells.Item(25, 2).Value = "NSE:NV20BEES"
$ws.Cells.Item(25, 3).Value = ""
$ws.Cells.Item(25, 4).Value = ""
$ws.Cells.Item(25, 5).Value = ""
$ws.Cells.Item(25, 6).Value = ""

# Row 26
$ws.Cells.Item(26, 2).Value = "NSE:ONELIFECAP"
$ws.Cells.Item(26, 3).Value = ""
$ws.Cells.Item(26, 4).Value = ""
$ws.Cells.Item(26, 5).Value = ""
$ws.Cells.Item(26, 6).Value = ""

# Row 27
$ws.Cells.Item(27, 2).Value = "NSE:ORCHPHARMA"
$ws.Cells.Item(27, 3).Value = ""
$ws.Cells.Item(27, 4).Value = ""
$ws.Cells.Item(27, 5).Value = ""
$ws.Cells.Item(27, 6).Value = ""

# Row 28
$ws.Cells.Item(28, 2).Value = "NSE:PANACHE"
$ws.Cells.Item(28, 3).Value = ""
$ws.Cells.Item(28, 4).Value = ""
$ws.Cells.Item(28, 5).Value = ""
$ws.Cells.Item(28, 6).Value = ""

# Row 29
$ws.Cells.Item(29, 2).Value = "NSE:PILANIINVS"
$ws.Cells.Item(29, 3).Value = ""
$ws.Cells.Item(29, 4).Value = ""
$ws.Cells.Item(29, 5).Value = ""
$ws.Cells.Item(29, 6).Value = ""

# Row 30
$ws.Cells.Item(30, 2).Value = "NSE:PRIVISCL"
$ws.Cells.Item(30, 3).Value = ""
$ws.Cells.Item(30, 4).Value = ""
$ws.Cells.Item(30, 5).Value = ""
$ws.Cells.Item(30, 6).Value = ""

# Row 31
$ws.Cells.Item(31, 2).Value = "NSE:RANASUG"
$ws.Cells.Item(31, 3).Value = ""
$ws.Cells.Item(31, 4).Value = ""
$ws.Cells.Item(31, 5).Value = ""
$ws.Cells.Item(31, 6).Value = ""

# Row 32
$ws.Cells.Item(32, 2).Value = "NSE:REPRO"
$ws.Cells.Item(32, 3).Value = ""
$ws.Cells.Item(32, 4).Value = ""
$ws.Cells.Item(32, 5).Value = ""
$ws.Cells.Item(32, 6).Value = ""
